$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 361; this shifts rows 361..433 down to 362..434
$ws.Rows.Item(361).Insert()

# Populate the newly inserted row 361 with a new weekly data point
# (same shape/values as the row that used to be at 361, except the date)
$ws.Cells.Item(361, 1).Value2  = 4
$ws.Cells.Item(361, 2).Value2  = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(361, 3).Value2  = "Los Lagos"
$ws.Cells.Item(361, 4).Value2  = 45015
$ws.Cells.Item(361, 5).Value2  = 10
$ws.Cells.Item(361, 6).Value2  = 100112040
$ws.Cells.Item(361, 7).Value2  = "Cilantro"
$ws.Cells.Item(361, 8).Value2  = "Sin especificar"
$ws.Cells.Item(361, 9).Value2  = "Primera"
$ws.Cells.Item(361, 10).Value2 = 50
$ws.Cells.Item(361, 11).Value2 = 6000
$ws.Cells.Item(361, 12).Value2 = 6000
$ws.Cells.Item(361, 13).Value2 = 6000
$ws.Cells.Item(361, 14).Value2 = "`$/docena de atados (2 kilos)"
$ws.Cells.Item(361, 15).Value2 = "Región de La Araucanía"
$ws.Cells.Item(361, 16).Value2 = 3000
$ws.Cells.Item(361, 17).Value2 = 2
$ws.Cells.Item(361, 18).Value2 = "Hortaliza"

# Make sure the new D cell keeps the date number format used by the rest of column D
$ws.Cells.Item(361, 4).NumberFormat = $ws.Cells.Item(362, 4).NumberFormat
